$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3334332.5
$ws.Range("I33").Value = 3334332.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 3334332.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -3334103.5
$ws.Range("N33").ClearContents()
$ws.Range("H62").Value = 4634993
$ws.Range("I62").Value = 6948490
$ws.Range("K62").Value = 6948490
$ws.Range("M62").Value = -6947866
$ws.Range("H65").Value = 4634993
$ws.Range("I65").Value = 6948490
$ws.Range("K65").Value = 34742450
$ws.Range("M65").Value = -34739330
$ws.Range("H86").Value = 1647168.9
$ws.Range("I86").Value = 1505.5555
$ws.Range("J86").Value = 3763021.8
$ws.Range("K86").Value = 1505.5555
$ws.Range("L86").Value = 3763021.8
$ws.Range("M86").Value = -382.5554999999999
$ws.Range("N86").Value = -3765267.8
$ws.Range("H89").Value = 1647168.9
$ws.Range("I89").Value = 1505.5555
$ws.Range("J89").Value = 3763021.8
$ws.Range("K89").Value = 7527.7775
$ws.Range("L89").Value = 18815109
$ws.Range("M89").Value = -1911.7775
$ws.Range("N89").Value = -18826341
$ws.Range("H97").Value = 3128.4614
$ws.Range("J97").Value = 3128.4614
$ws.Range("L97").Value = 9385.3842
$ws.Range("N97").Value = -10377.3842
$ws.Range("H100").Value = 10753.056
$ws.Range("I100").Value = 4226.4287
$ws.Range("J100").Value = 14906.363
$ws.Range("K100").Value = 4226.4287
$ws.Range("L100").Value = 14906.363
$ws.Range("M100").Value = -3685.4287
$ws.Range("N100").Value = -15988.363
$ws.Range("H107").Value = 114667.445
$ws.Range("I107").Value = 169668
$ws.Range("J107").Value = 4666.3335
$ws.Range("K107").Value = 169668
$ws.Range("L107").Value = 4666.3335
$ws.Range("M107").Value = -167748
$ws.Range("N107").Value = -8506.333500000001
$ws.Range("H132").Value = 2127.3953
$ws.Range("I132").Value = 2073.8206
$ws.Range("K132").Value = 6221.4618
$ws.Range("M132").Value = -3691.4618

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5470.902
$ws.Range("I32").Value = 4681.2554
$ws.Range("K32").Value = 4681.2554
$ws.Range("M32").Value = -4394.2554
$ws.Range("H61").Value = 3650.2646
$ws.Range("I61").Value = 3537.3462
$ws.Range("K61").Value = 3537.3462
$ws.Range("M61").Value = -3325.3462
$ws.Range("H74").Value = 3642.111
$ws.Range("I74").Value = 6970.3335
$ws.Range("J74").Value = 1978
$ws.Range("K74").Value = 6970.3335
$ws.Range("L74").Value = 1978
$ws.Range("M74").Value = -6096.3335
$ws.Range("N74").Value = -3726
$ws.Range("H77").Value = 3642.111
$ws.Range("I77").Value = 6970.3335
$ws.Range("J77").Value = 1978
$ws.Range("K77").Value = 34851.6675
$ws.Range("L77").Value = 9890
$ws.Range("M77").Value = -30483.6675
$ws.Range("N77").Value = -18626
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 4093.3394
$ws.Range("J132").Value = 4880.4346
$ws.Range("L132").Value = 14641.3038
$ws.Range("N132").Value = -19701.3038
$ws.Range("H136").Value = 3650.2646
$ws.Range("I136").Value = 3537.3462
$ws.Range("K136").Value = 10612.0386
$ws.Range("M136").Value = -8062.0386

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 155.72728
$ws.Range("I80").Value = 117.75
$ws.Range("K80").Value = 117.75
$ws.Range("M80").Value = 880.25
$ws.Range("H83").Value = 155.72728
$ws.Range("I83").Value = 117.75
$ws.Range("K83").Value = 588.75
$ws.Range("M83").Value = 4403.25
$ws.Range("H99").Value = 3689.8
$ws.Range("I99").Value = 3493.5
$ws.Range("K99").Value = 3493.5
$ws.Range("M99").Value = -1995.5
$ws.Range("H107").Value = 334694.8
$ws.Range("I107").Value = 1162.9546
$ws.Range("J107").Value = 1251907.4
$ws.Range("K107").Value = 1162.9546
$ws.Range("L107").Value = 1251907.4
$ws.Range("M107").Value = 757.0454
$ws.Range("N107").Value = -1255747.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4266.926
$ws.Range("I16").Value = 3815.158
$ws.Range("J16").Value = 5339.875
$ws.Range("K16").Value = 3815.158
$ws.Range("L16").Value = 5339.875
$ws.Range("M16").Value = -3528.158
$ws.Range("N16").Value = -5913.875
$ws.Range("H31").Value = 36761.31
$ws.Range("I31").Value = 1078.3846
$ws.Range("J31").Value = 65753.69
$ws.Range("K31").Value = 1078.3846
$ws.Range("L31").Value = 65753.69
$ws.Range("M31").Value = -783.3846000000001
$ws.Range("N31").Value = -66343.69
$ws.Range("H34").Value = 36761.31
$ws.Range("I34").Value = 1078.3846
$ws.Range("J34").Value = 65753.69
$ws.Range("K34").Value = 1078.3846
$ws.Range("L34").Value = 65753.69
$ws.Range("M34").Value = -876.3846000000001
$ws.Range("N34").Value = -66157.69
$ws.Range("H105").Value = 1067.5
$ws.Range("I105").Value = 1067.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1067.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 679.5
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 4266.926
$ws.Range("I113").Value = 3815.158
$ws.Range("J113").Value = 5339.875
$ws.Range("K113").Value = 3815.158
$ws.Range("L113").Value = 5339.875
$ws.Range("M113").Value = -1645.158
$ws.Range("N113").Value = -9679.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 569234.5600000001
$ws.Range("I5").Value = 47685.293
$ws.Range("J5").Value = 5002403.5
$ws.Range("K5").Value = 143055.879
$ws.Range("L5").Value = 15007210.5
$ws.Range("M5").Value = -142943.879
$ws.Range("N5").Value = -15007434.5
$ws.Range("H135").Value = 569234.5600000001
$ws.Range("I135").Value = 47685.293
$ws.Range("J135").Value = 5002403.5
$ws.Range("K135").Value = 429167.637
$ws.Range("L135").Value = 45021631.5
$ws.Range("M135").Value = -426632.637
$ws.Range("N135").Value = -45026701.5
$ws.Range("H139").Value = 5253.278
$ws.Range("I139").Value = 1456.6111
$ws.Range("K139").Value = 4369.8333
$ws.Range("M139").Value = 770.1666999999998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 676763.3
$ws.Range("I113").Value = 1431142.9
$ws.Range("K113").Value = 1431142.9
$ws.Range("M113").Value = -1428972.9
$ws.Range("H139").Value = 89333.336
$ws.Range("J139").Value = 89333.336
$ws.Range("L139").Value = 89333.336
$ws.Range("N139").Value = -99613.336

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 479.125
$ws.Range("I22").Value = 581.6667
$ws.Range("J22").Value = 417.6
$ws.Range("K22").Value = 581.6667
$ws.Range("L22").Value = 417.6
$ws.Range("M22").Value = -286.6667
$ws.Range("N22").Value = -1007.6
$ws.Range("H27").Value = 479.125
$ws.Range("I27").Value = 581.6667
$ws.Range("J27").Value = 417.6
$ws.Range("K27").Value = 581.6667
$ws.Range("L27").Value = 417.6
$ws.Range("M27").Value = -474.6667
$ws.Range("N27").Value = -631.6
$ws.Range("H132").Value = 2757.4
$ws.Range("I132").Value = 2403.4546
$ws.Range("K132").Value = 7210.3638
$ws.Range("M132").Value = -4680.3638

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 1500
$ws.Range("J25").Value = 1500
$ws.Range("L25").Value = 1500
$ws.Range("N25").Value = -2086
$ws.Range("H113").Value = 757.0714
$ws.Range("I113").Value = 742.8261
$ws.Range("K113").Value = 2228.4783
$ws.Range("M113").Value = -58.47829999999976
$ws.Range("H122").Value = 27031008
$ws.Range("I122").Value = 37040956
$ws.Range("K122").Value = 111122868
$ws.Range("M122").Value = -111120418
$ws.Range("H126").Value = 1533
$ws.Range("I126").Value = 1533
$ws.Range("K126").Value = 4599
$ws.Range("M126").Value = -2129

